# Generate Report for Handoff
#
# The a378deed-... file has now been handed off (status "Ready for handoff"),
# while b673f8c4-... remains "Handed back: in sync with en-US".  Since rows
# are kept sorted (most-recently-updated / pending items first), the two
# data rows on every sheet swap places: row 2 now holds the b673f8c4 data
# and row 3 now holds the a378deed data (previously the other way round).
# The a378deed row also gets its Status / timestamp fields refreshed.

$wb = $excel.ActiveWorkbook

function Swap-HyperlinkDisplays($ws, $pairs) {
    $hls = @()
    foreach ($hl in $ws.Hyperlinks) {
        $hls += $hl
    }
    $map = @{}
    foreach ($hl in $hls) {
        $addr = $hl.Range.Address(0, 0)
        $map[$addr] = $hl
    }
    foreach ($p in $pairs) {
        $addr1 = $p[0]
        $addr2 = $p[1]
        $h1 = $map[$addr1]
        $h2 = $map[$addr2]
        $t1 = $h1.TextToDisplay
        $t2 = $h2.TextToDisplay
        $h1.TextToDisplay = $t2
        $h2.TextToDisplay = $t1
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2 becomes the b673f8c4 entry, row 3 becomes the a378deed entry.
$wsOverview.Range("A2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.md"
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("D2").Value = "2016-03-22 02:40:58"

$wsOverview.Range("A3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 02:42:18"

Swap-HyperlinkDisplays $wsOverview @(, @("A2", "A3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.98f4df468d91f75d64cf0695095d70b3fafb3b52.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-22 02:40:54"
$wsZh.Range("F2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.md"
$wsZh.Range("G2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.98f4df468d91f75d64cf0695095d70b3fafb3b52.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-03-22 02:41:30"
$wsZh.Range("J2").Value = "Include"

$wsZh.Range("A3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.7bc9057f83e78f2e26832acdbca5d8a7adffff90.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-22 02:42:14"
$wsZh.Range("F3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.md"
$wsZh.Range("G3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.7bc9057f83e78f2e26832acdbca5d8a7adffff90.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-03-22 02:41:30"
$wsZh.Range("J3").Value = "Include"

Swap-HyperlinkDisplays $wsZh @(@("A2", "A3"), @("D2", "D3"), @("F2", "F3"), @("G2", "G3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.98f4df468d91f75d64cf0695095d70b3fafb3b52.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-22 02:40:58"
$wsDe.Range("F2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.md"
$wsDe.Range("G2").Value = "b673f8c4-9f47-49de-aa31-d6ae0ad5158f.98f4df468d91f75d64cf0695095d70b3fafb3b52.de-de.xlf"
$wsDe.Range("H2").Value = "2016-03-22 02:41:37"
$wsDe.Range("J2").Value = "Include"

$wsDe.Range("A3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.7bc9057f83e78f2e26832acdbca5d8a7adffff90.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-22 02:40:58"
$wsDe.Range("F3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.md"
$wsDe.Range("G3").Value = "a378deed-3832-42ed-b959-49ca0c0bd1a1.7bc9057f83e78f2e26832acdbca5d8a7adffff90.de-de.xlf"
$wsDe.Range("H3").Value = "2016-03-22 02:41:37"
$wsDe.Range("J3").Value = "Include"

Swap-HyperlinkDisplays $wsDe @(@("A2", "A3"), @("D2", "D3"), @("F2", "F3"), @("G2", "G3"))
